$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "47.050.14"
$ws.Range("E2").Value = "  +5.40%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.352.02"
$ws.Range("E3").Value = "  +4.75%  "

$ws.Range("E4").Value = "  -0.78%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.95"
$ws.Range("E5").Value = "  -0.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.14"
$ws.Range("E6").Value = "  +3.50%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.578"
$ws.Range("E7").Value = "  +1.44%  "

$ws.Range("E8").Value = "  -0.69%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.538"
$ws.Range("E9").Value = "  +4.38%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.92"
$ws.Range("E10").Value = "  +2.78%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0808"
$ws.Range("E11").Value = "  +0.84%  "

$ws.Range("E12").Value = "  +3.31%  "

$ws.Range("E13").Value = "  -0.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.707.51"
$ws.Range("E14").Value = "  +4.53%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.349.44"
$ws.Range("E15").Value = "  +4.86%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.26"
$ws.Range("E16").Value = "  +5.08%  "

$ws.Range("E17").Value = "  +0.00%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "46.859.57"
$ws.Range("E18").Value = "  +5.47%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.71"
$ws.Range("E19").Value = "  +17.00%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0951"
$ws.Range("E20").Value = "  +1.56%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.20"
$ws.Range("E21").Value = "  +0.39%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.87"
$ws.Range("E22").Value = "  +2.36%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "245.69"
$ws.Range("E23").Value = "  +3.59%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.98"
$ws.Range("E24").Value = "  +0.94%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.99"
$ws.Range("E25").Value = "  +1.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.47%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "41.97"
$ws.Range("E27").Value = "  +13.44%  "

$ws.Range("E28").Value = "  -0.32%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.90"
$ws.Range("E29").Value = "  +1.31%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.21"
$ws.Range("E30").Value = "  +1.13%  "

$ws.Range("E31").Value = "  -1.64%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "152.64"
$ws.Range("E32").Value = "  +3.69%  "

$ws.Range("E33").Value = "  +4.54%  "

$ws.Range("E34").Value = "  +0.37%  "

$ws.Range("E35").Value = "  -0.63%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.110"
$ws.Range("E36").Value = "  +1.69%  "

$ws.Range("E37").Value = "  +0.80%  "

$ws.Range("E38").Value = "  -0.88%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.06"
$ws.Range("E39").Value = "  +7.53%  "

$ws.Range("E40").Value = "  +5.35%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.42"
$ws.Range("E41").Value = "  +2.32%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.86"
$ws.Range("E42").Value = "  -8.75%  "

$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.01"
$ws.Range("E43").Value = "  +13.74%  "

$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  -0.84%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.852.77"
$ws.Range("E45").Value = "  +2.22%  "

$ws.Range("E46").Value = "  +4.96%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "74.46"
$ws.Range("E47").Value = "  +8.07%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "81.21"
$ws.Range("E48").Value = "  -0.78%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.95"
$ws.Range("E49").Value = "  +2.83%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "98.70"
$ws.Range("E50").Value = "  +0.47%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "55.39"
$ws.Range("E51").Value = "  +2.47%  "
